# Auto-generated edit script: updates market-price-derived columns (H:N)
# on multiple worksheets, matching a scheduled market-data refresh.
$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H33").Value = 691.8823
$ws.Range("I33").Value = 626.2857
$ws.Range("K33").Value = 626.2857
$ws.Range("M33").Value = -397.2857
$ws.Range("H54").Value = 22333.334
$ws.Range("I54").Value = 16000
$ws.Range("J54").Value = 35000
$ws.Range("K54").Value = 16000
$ws.Range("L54").Value = 35000
$ws.Range("M54").Value = -15514
$ws.Range("N54").Value = -35972
$ws.Range("H113").Value = 2875
$ws.Range("I113").Value = 2772.6667
$ws.Range("J113").Value = 3489
$ws.Range("K113").Value = 2772.6667
$ws.Range("L113").Value = 3489
$ws.Range("M113").Value = 481.3332999999998
$ws.Range("N113").Value = -9997
$ws.Range("H131").Value = 2898.2856
$ws.Range("I131").Value = 2898.2856
$ws.Range("J131").Value = 0
$ws.Range("K131").Value = 8694.856800000001
$ws.Range("L131").Value = 0
$ws.Range("M131").Value = -3654.856800000001
$ws.Range("N131").ClearContents()
$ws.Range("H132").Value = 266113.72
$ws.Range("I132").Value = 320212.03
$ws.Range("J132").Value = 7643.8887
$ws.Range("K132").Value = 960636.0900000001
$ws.Range("L132").Value = 22931.6661
$ws.Range("M132").Value = -958106.0900000001
$ws.Range("N132").Value = -27991.6661
$ws.Range("H133").Value = 81100
$ws.Range("J133").Value = 81100
$ws.Range("L133").Value = 81100
$ws.Range("N133").Value = -91220
$ws.Range("H138").Value = 3119.6836
$ws.Range("J138").Value = 3793.8948
$ws.Range("L138").Value = 11381.6844
$ws.Range("N138").Value = -21661.6844
$ws.Range("H139").Value = 254135
$ws.Range("J139").Value = 254135
$ws.Range("L139").Value = 254135
$ws.Range("N139").Value = -264415
$ws.Range("H141").Value = 3207.2222
$ws.Range("I141").Value = 2996.8572
$ws.Range("K141").Value = 8990.571599999999
$ws.Range("M141").Value = -3810.571599999999

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 1889047.2
$ws.Range("I32").Value = 4418.24
$ws.Range("K32").Value = 4418.24
$ws.Range("M32").Value = -4131.24
$ws.Range("H74").Value = 4909.1724
$ws.Range("I74").Value = 4598.6665
$ws.Range("K74").Value = 4598.6665
$ws.Range("M74").Value = -3724.6665
$ws.Range("H77").Value = 4909.1724
$ws.Range("I77").Value = 4598.6665
$ws.Range("K77").Value = 22993.3325
$ws.Range("M77").Value = -18625.3325
$ws.Range("H97").Value = 2053.1155
$ws.Range("I97").Value = 1654.2632
$ws.Range("J97").Value = 3135.7144
$ws.Range("K97").Value = 1654.2632
$ws.Range("L97").Value = 3135.7144
$ws.Range("M97").Value = -1158.2632
$ws.Range("N97").Value = -4127.7144
$ws.Range("H122").Value = 4227.4194
$ws.Range("I122").Value = 3205.6296
$ws.Range("K122").Value = 9616.888800000001
$ws.Range("M122").Value = -7166.888800000001
$ws.Range("H139").Value = 77406.82000000001
$ws.Range("J139").Value = 77406.82000000001
$ws.Range("L139").Value = 77406.82000000001
$ws.Range("N139").Value = -87686.82000000001

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 4732.731
$ws.Range("J86").Value = 7989.5454
$ws.Range("L86").Value = 7989.5454
$ws.Range("N86").Value = -10235.5454
$ws.Range("H89").Value = 4732.731
$ws.Range("J89").Value = 7989.5454
$ws.Range("L89").Value = 39947.727
$ws.Range("N89").Value = -51179.727
$ws.Range("H94").Value = 4076.068
$ws.Range("I94").Value = 2948.0334
$ws.Range("J94").Value = 6493.2856
$ws.Range("K94").Value = 2948.0334
$ws.Range("L94").Value = 6493.2856
$ws.Range("M94").Value = -2497.0334
$ws.Range("N94").Value = -7395.2856
$ws.Range("H107").Value = 3852573.2
$ws.Range("I107").Value = 4767662
$ws.Range("K107").Value = 4767662
$ws.Range("M107").Value = -4765742

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 9808.297
$ws.Range("I31").Value = 16259
$ws.Range("J31").Value = 6314.1665
$ws.Range("K31").Value = 16259
$ws.Range("L31").Value = 6314.1665
$ws.Range("M31").Value = -15964
$ws.Range("N31").Value = -6904.1665
$ws.Range("H34").Value = 9808.297
$ws.Range("I34").Value = 16259
$ws.Range("J34").Value = 6314.1665
$ws.Range("K34").Value = 16259
$ws.Range("L34").Value = 6314.1665
$ws.Range("M34").Value = -16057
$ws.Range("N34").Value = -6718.1665
$ws.Range("H94").Value = 5950.5
$ws.Range("I94").Value = 2599
$ws.Range("K94").Value = 2599
$ws.Range("M94").Value = -2148
$ws.Range("H132").Value = 8162.2856
$ws.Range("J132").Value = 18134.857
$ws.Range("L132").Value = 54404.571
$ws.Range("N132").Value = -59464.571
$ws.Range("H133").Value = 0
$ws.Range("J133").Value = 0
$ws.Range("L133").Value = 0
$ws.Range("N133").ClearContents()
$ws.Range("H141").Value = 231707.08
$ws.Range("J141").Value = 247266
$ws.Range("L141").Value = 247266
$ws.Range("N141").Value = -257626

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H17").Value = 2762.1667
$ws.Range("J17").Value = 1616.1666
$ws.Range("L17").Value = 4848.4998
$ws.Range("N17").Value = -5186.4998
$ws.Range("H42").Value = 3000
$ws.Range("I42").Value = 3000
$ws.Range("K42").Value = 9000
$ws.Range("M42").Value = -8466
$ws.Range("H62").Value = 19666.666
$ws.Range("I62").Value = 18000
$ws.Range("K62").Value = 54000
$ws.Range("M62").Value = -53314
$ws.Range("H65").Value = 19666.666
$ws.Range("I65").Value = 18000
$ws.Range("K65").Value = 162000
$ws.Range("M65").Value = -158568
$ws.Range("H121").Value = 39774.875
$ws.Range("I121").Value = 2962.75
$ws.Range("J121").Value = 76587
$ws.Range("K121").Value = 8888.25
$ws.Range("L121").Value = 229761
$ws.Range("M121").Value = -7578.25
$ws.Range("N121").Value = -232381
$ws.Range("H131").Value = 1794.8823
$ws.Range("I131").Value = 1364.4166
$ws.Range("J131").Value = 2828
$ws.Range("K131").Value = 4093.2498
$ws.Range("L131").Value = 8484
$ws.Range("M131").Value = 946.7501999999999
$ws.Range("N131").Value = -18564

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 4888.768
$ws.Range("I102").Value = 3602.5557
$ws.Range("J102").Value = 7203.95
$ws.Range("K102").Value = 3602.5557
$ws.Range("L102").Value = 7203.95
$ws.Range("M102").Value = -1980.5557
$ws.Range("N102").Value = -10447.95

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H93").Value = 4754.4165
$ws.Range("I93").Value = 2885.875
$ws.Range("J93").Value = 8491.5
$ws.Range("K93").Value = 2885.875
$ws.Range("L93").Value = 8491.5
$ws.Range("M93").Value = -1637.875
$ws.Range("N93").Value = -10987.5
$ws.Range("H100").Value = 3638.9565
$ws.Range("I100").Value = 6799.2
$ws.Range("K100").Value = 6799.2
$ws.Range("M100").Value = -6258.2

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H95").Value = 36999.668
$ws.Range("J95").Value = 36999.668
$ws.Range("L95").Value = 36999.668
$ws.Range("N95").Value = -42491.668
$ws.Range("H113").Value = 4386986
$ws.Range("I113").Value = 7576752.5
$ws.Range("J113").Value = 1056.5625
$ws.Range("K113").Value = 22730257.5
$ws.Range("L113").Value = 3169.6875
$ws.Range("M113").Value = -22728087.5
$ws.Range("N113").Value = -7509.6875
$ws.Range("H135").Value = 63932.8
$ws.Range("J135").Value = 63932.8
$ws.Range("L135").Value = 63932.8
$ws.Range("N135").Value = -74072.8
$ws.Range("H136").Value = 7691.25
$ws.Range("I136").Value = 7406.0977
$ws.Range("K136").Value = 22218.2931
$ws.Range("M136").Value = -19668.2931
